# Update the "想去人数" (want-to-go headcount) figures in column F
# for the "展览" and "全部类型" worksheets, reflecting a refreshed
# data pull (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------
$wsExhibition = $wb.Worksheets.Item("展览")

$exhibitionUpdates = @{
    2  = 37
    4  = 16288
    5  = 423
    6  = 17
    7  = 732
    8  = 15582
    9  = 66
    10 = 9227
    11 = 462
    13 = 1027
    14 = 120
    15 = 214
    17 = 219
    19 = 86
    20 = 602
    23 = 74
    24 = 1149
    26 = 20
    28 = 520
    34 = 62
    36 = 364
    37 = 476
    38 = 120
    39 = 5664
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# --- Sheet "全部类型" -------------------------------------------------
$wsAllTypes = $wb.Worksheets.Item("全部类型")

$allTypesUpdates = @{
    2  = 37
    4  = 16289
    5  = 423
    6  = 17
    7  = 732
    8  = 15582
    9  = 66
    10 = 9227
    11 = 462
    13 = 1027
    14 = 120
    15 = 214
    17 = 219
    19 = 86
    20 = 602
    23 = 74
    24 = 1149
    26 = 20
    28 = 520
    36 = 62
    38 = 364
    39 = 476
    40 = 120
    41 = 5664
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
